{"js": "// Apply the \"Lastenheft\" text revisions described by the commit.\n//\n// 1) \"Greenfoot Game \u2026 . Das Spiel ...\" -> \"Greenfoot-Game. Das Spiel ...\"\n// 2) \"Das Ziel des Spiels ist es alle Welten\" -> \"... ist es, alle Welten\"\n// 3) \"ihre St\u00e4rke und Schnelligkeit, die Waffen\" -> \"... Schnelligkeit, Waffen\"\n// 4) \"vorrankommt\" -> \"vorankommt\"\n// 5) Merge the stray empty paragraph into the following \"Durch ein\n//    Lebenssystem...\" paragraph, rework its wording and fix \"an oberen\"\n//    -> \"am oberen\".\n\n// 1) \"Greenfoot Game \u2026 .\" -> \"Greenfoot-Game.\"\nlet results = context.document.body.search(\" Game \u2026 . \", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"-Game. \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Word stamps the cursor's last position with a \"_GoBack\" bookmark on\n// save; reproduce it right after the second \"Greenfoot\" mention.\nresults = context.document.body.search(\"Greenfoot\", { matchCase: true });\nawait context.sync();\nif (results.items.length > 1) {\n  results.items[1].getRange(\"End\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) add comma: \"ist es alle Welten\" -> \"ist es, alle Welten\"\nresults = context.document.body.search(\"ist es alle Welten\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"ist es, alle Welten\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) drop \"die \" before \"Waffen\"\nresults = context.document.body.search(\"Schnelligkeit, die Waffen\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"Schnelligkeit, Waffen\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) typo fix: \"vorrankommt\" -> \"vorankommt\"\nresults = context.document.body.search(\"vorrankommt\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"vorankommt\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 5a) remove the empty paragraph between \"...bewegt.\" and \"Durch ein\n//     Lebenssystem...\", merging it into the next paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\") {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 5b) move \", wenn er von\" right after \"Durch ein\" and drop it from its\n//     old spot after \"sterben\".\nresults = context.document.body.search(\n  \"Durch ein Lebenssystem kann der Spieler auch sterben, wenn er von Gegnern\",\n  { matchCase: true }\n);\nawait context.sync();\nresults.items[0].insertText(\n  \"Durch ein, wenn er von Lebenssystem kann der Spieler auch sterben Gegnern\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 5c) \"Lebensanzeige an oberen\" -> \"Lebensanzeige am oberen\"\nresults = context.document.body.search(\"Lebensanzeige an oberen\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"Lebensanzeige am oberen\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply the \"Lastenheft\" text revisions described by the commit.\n#\n# 1) \"Greenfoot Game \u2026 . Das Spiel ...\" -> \"Greenfoot-Game. Das Spiel ...\"\n# 2) \"Das Ziel des Spiels ist es alle Welten\" -> \"... ist es, alle Welten\"\n# 3) \"ihre Staerke und Schnelligkeit, die Waffen\" -> \"... Schnelligkeit, Waffen\"\n# 4) \"vorrankommt\" -> \"vorankommt\"\n# 5) Merge the stray empty paragraph into the following \"Durch ein\n#    Lebenssystem...\" paragraph, rework its wording and fix \"an oberen\"\n#    -> \"am oberen\".\n\n$d = $word.ActiveDocument\n\n# 1) \"Greenfoot Game \u2026 .\" -> \"Greenfoot-Game.\"\n$find = $d.Content.Find\n$find.Execute(\" Game \u2026 . \", $false, $false, $false, $false, $false, $true, 1, $false, \"-Game. \", 2)\n\n# Word stamps the cursor's last edit position with a \"_GoBack\" bookmark on\n# save; reproduce it right after the second \"Greenfoot\" mention.\n$r1 = $d.Content\n$r1.Find.Execute(\"Greenfoot\")\n$r2 = $d.Range($r1.End, $d.Content.End)\n$r2.Find.Execute(\"Greenfoot\")\n$bmRange = $d.Range($r2.End, $r2.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# 2) add comma: \"ist es alle Welten\" -> \"ist es, alle Welten\"\n$find2 = $d.Content.Find\n$find2.Execute(\"ist es alle Welten\", $false, $false, $false, $false, $false, $true, 1, $false, \"ist es, alle Welten\", 2)\n\n# 3) drop \"die \" before \"Waffen\"\n$find3 = $d.Content.Find\n$find3.Execute(\"Schnelligkeit, die Waffen\", $false, $false, $false, $false, $false, $true, 1, $false, \"Schnelligkeit, Waffen\", 2)\n\n# 4) typo fix: \"vorrankommt\" -> \"vorankommt\"\n$find4 = $d.Content.Find\n$find4.Execute(\"vorrankommt\", $false, $false, $false, $false, $false, $true, 1, $false, \"vorankommt\", 2)\n\n# 5a) remove the empty paragraph between \"...bewegt.\" and \"Durch ein\n#     Lebenssystem...\", merging it into the next paragraph.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim().Length -eq 0) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 5b) move \", wenn er von\" right after \"Durch ein\" and drop it from its\n#     old spot after \"sterben\".\n$find5 = $d.Content.Find\n$find5.Execute(\"Durch ein Lebenssystem kann der Spieler auch sterben, wenn er von Gegnern\", $false, $false, $false, $false, $false, $true, 1, $false, \"Durch ein, wenn er von Lebenssystem kann der Spieler auch sterben Gegnern\", 2)\n\n# 5c) \"Lebensanzeige an oberen\" -> \"Lebensanzeige am oberen\"\n$find6 = $d.Content.Find\n$find6.Execute(\"Lebensanzeige an oberen\", $false, $false, $false, $false, $false, $true, 1, $false, \"Lebensanzeige am oberen\", 2)\n"}
